# Applies the Wed May 22 08:00:15 UTC 2024 cryptos list refresh (GitHub Actions).
# Updates Price (col D) / Volume(1h) (col E) values, and swaps the
# dogwifhat / Bittensor rows (43 <-> 44) with their refreshed figures.
# Numeric-looking Price strings are apostrophe-prefixed so Excel keeps
# them as text (matching the original inline-string cell type) instead
# of auto-converting them to numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '69.801.77'
$ws.Range("E2").Value = '  -1.61%  '
$ws.Range("D3").Value = '3.749.32'
$ws.Range("E3").Value = '  +2.36%  '
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").Value = '''622.04'
$ws.Range("E5").Value = '  -1.20%  '
$ws.Range("D6").Value = '''180.18'
$ws.Range("E6").Value = '  -0.55%  '
$ws.Range("D7").Value = '3.747.92'
$ws.Range("E7").Value = '  +2.46%  '
$ws.Range("D8").Value = '''0.999'
$ws.Range("E8").Value = '  -0.08%  '
$ws.Range("D9").Value = '''0.534'
$ws.Range("E9").Value = '  -1.02%  '
$ws.Range("E10").Value = '  +3.09%  '
$ws.Range("E11").Value = '  -5.60%  '
$ws.Range("E12").Value = '  -2.16%  '
$ws.Range("D13").Value = '''41.01'
$ws.Range("E13").Value = '  +0.94%  '
$ws.Range("D14").Value = '''0.0000262'
$ws.Range("E14").Value = '  +2.98%  '
$ws.Range("D15").Value = '4.367.62'
$ws.Range("D16").Value = '3.749.82'
$ws.Range("E16").Value = '  +1.93%  '
$ws.Range("D17").Value = '69.832.61'
$ws.Range("E17").Value = '  -1.59%  '
$ws.Range("D18").Value = '''0.123'
$ws.Range("E18").Value = '  -0.07%  '
$ws.Range("E19").Value = '  +1.47%  '
$ws.Range("D20").Value = '''16.72'
$ws.Range("E20").Value = '  -1.18%  '
$ws.Range("D21").Value = '''506.04'
$ws.Range("E21").Value = '  -2.60%  '
$ws.Range("D22").Value = '''9.49'
$ws.Range("E22").Value = '  +2.25%  '
$ws.Range("E23").Value = '  -1.89%  '
$ws.Range("E24").Value = '  +0.83%  '
$ws.Range("D25").Value = '''87.20'
$ws.Range("E25").Value = '  -1.45%  '
$ws.Range("D26").Value = '''13.19'
$ws.Range("E26").Value = '  -2.12%  '
$ws.Range("D27").Value = '''11.14'
$ws.Range("E27").Value = '  +1.39%  '
$ws.Range("E28").Value = '  +25.76%  '
$ws.Range("E30").Value = '  -1.38%  '
$ws.Range("D31").Value = '''2.93'
$ws.Range("E31").Value = '  -0.88%  '
$ws.Range("E32").Value = '  -2.71%  '
$ws.Range("D33").Value = '''31.20'
$ws.Range("E33").Value = '  -1.29%  '
$ws.Range("D34").Value = '''0.116'
$ws.Range("E34").Value = '  -0.02%  '
$ws.Range("E35").Value = '  -0.08%  '
$ws.Range("E36").Value = '  +3.71%  '
$ws.Range("D37").Value = '''6.21'
$ws.Range("E37").Value = '  +1.05%  '
$ws.Range("E38").Value = '  -3.22%  '
$ws.Range("E39").Value = '  +2.00%  '
$ws.Range("E40").Value = '  -4.78%  '
$ws.Range("D41").Value = '''49.96'
$ws.Range("E41").Value = '  -3.65%  '
$ws.Range("D42").Value = '''45.68'
$ws.Range("E42").Value = '  +0.31%  '
$ws.Range("B43").Value = 'dogwifhat'
$ws.Range("C43").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D43").Value = '''2.92'
$ws.Range("E43").Value = '  +5.19%  '
$ws.Range("B44").Value = 'Bittensor'
$ws.Range("C44").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D44").Value = '''425.54'
$ws.Range("E44").Value = '  +0.30%  '
$ws.Range("E45").Value = '  -1.34%  '
$ws.Range("D46").Value = '3.008.15'
$ws.Range("E46").Value = '  -3.68%  '
$ws.Range("D48").Value = '''27.41'
$ws.Range("E48").Value = '  -3.81%  '
$ws.Range("E49").Value = '  -0.04%  '
$ws.Range("D50").Value = '''138.04'
$ws.Range("E50").Value = '  -1.07%  '
$ws.Range("D51").Value = '''2.50'
$ws.Range("E51").Value = '  +1.50%  '